$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix a few existing values in rows 14, 18, 29, 31 (days already entered)
# ---------------------------------------------------------------------------
$ws.Range("G14").Value = 1
$ws.Range("J18").Value = 1
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 1
$ws.Range("J29").Value = 1
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 1

# ---------------------------------------------------------------------------
# 2) Append 35 new rows of data (rows 32-66) for day 6 (and some for day 7)
# ---------------------------------------------------------------------------
$aVals = @(45023,45023,45023,45023,45023,45023,45023,45023,45023,45023,45023,45023,45023,45023,45023,45023,45023,45023,45023,45023,45023,45023,45024,45024,45024,45024,45024,45024,45024,45024,45024,45024,45024,45024,45023)
$bVals = @(6296,2193,3107,5575,2383,9017,7330,2716,3089,8075,1330,751,4021,8510,4074,1082,9378,7607,9186,4612,9583,3389,5279,3022,3329,8465,3755,8779,1174,8050,9167,1997,2444,4326,7679)
$cVals = @("Resonators","Madolche","Dark Magician","Resonators","Phantom Knights","Timelord","Dark Magician","Starry Knight","Code Talker","Madolche","Resonators","Timelord","Phantom Knights","Odd-Eyes","Elemental Heros","Phantom Knights","Solfachord","Elemental Heros","Phantom Knights","D/D/D","Block Dragon","Timelord","Synchrons","Rose Dragon","Synchrons","Abyss Actor","Starry Knight","Resonators","Resonators","Salamangreat","Elemental Heros","Gouki","Elemental Heros","Gouki","Fire King")
$dVals = @("Alma Rugiente","Cambios","Magia del Mago","Alma Rugiente","Fantasma de la Traición","Vacío, Infinito y Luz Infinita","Magia del Mago","Sonrisa del Hada","Codificador Hablador Vivo","Robar sentido Nivel bajo","Alma Rugiente","Vacío, Infinito y Luz Infinita","Fantasma de la Traición","Xiangke y Xiangsheng","Alianza de Héroes","Fantasma de la Traición","Péndulos Unidos","Alianza de Héroes","Fantasma de la Traición","El Mando del Rey de la Perdición","Mundo Sonriente","Vacío, Infinito y Luz Infinita","Un Vínculo Ilumina el Futuro","Conjuro de Rosas","Un Vínculo Ilumina el Futuro","Tumbas Selladas","Sonrisa del Hada","Alma Rugiente","Alma Rugiente","Alzarse desde el Valle de las Llamas","Alianza de Héroes","El Evento Principal Gouki","Alianza de Héroes","El Evento Principal Gouki","Robo del Destino")
$fVals = @(1,1,1,1,1,1,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
$gVals = @(0,0,0,1,0,0,0,0,0,0,0,0,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
$hVals = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0)
$iVals = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,1,1,1,1,1,1,1,1,0)
$jVals = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,0,0,0,0,0,0,0,1)

$startRow = 32
$count = 35

# Set column B (Usuario) values first, before copying formats, so that the
# Text-formatted column does not coerce the numbers to strings.
for ($i = 0; $i -lt $count; $i++) {
  $r = $startRow + $i
  $ws.Cells.Item($r, 2).Value2 = $bVals[$i]
}

# Copy the formatting of the last existing row down across all the new rows.
$ws.Range("A31:K31").Copy()
$endRow = $startRow + $count - 1
$ws.Range("A" + $startRow + ":K" + $endRow).PasteSpecial(-4122)

# Fill in the remaining plain values.
for ($i = 0; $i -lt $count; $i++) {
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $aVals[$i]
  $ws.Cells.Item($r, 3).Value = $cVals[$i]
  $ws.Cells.Item($r, 4).Value = $dVals[$i]
  $ws.Cells.Item($r, 6).Value = $fVals[$i]
  $ws.Cells.Item($r, 7).Value = $gVals[$i]
  $ws.Cells.Item($r, 8).Value = $hVals[$i]
  $ws.Cells.Item($r, 9).Value = $iVals[$i]
  $ws.Cells.Item($r, 10).Value = $jVals[$i]
}

# Fill in the formulas for columns E and K (each one individually, matching
# the non-shared formula style used for these new rows).
for ($i = 0; $i -lt $count; $i++) {
  $r = $startRow + $i
  $eFormula = '=IF($A' + $r + '="","",$A' + $r + '-1)'
  $kFormula = '=IFERROR(ROUND((VALUE(TEXT($E' + $r + ',"DD"))),0),"")'
  $ws.Cells.Item($r, 5).Formula = $eFormula
  $ws.Cells.Item($r, 11).Formula = $kFormula
}

# ---------------------------------------------------------------------------
# 3) Convert the existing per-row formulas in E9:E31 / K9:K31 into shared
#    formulas (as Excel does when a formula is filled down across a range).
# ---------------------------------------------------------------------------
$ws.Range("E9:E31").Formula = '=IF($A9="","",$A9-1)'
$ws.Range("K9:K31").Formula = '=IFERROR(ROUND((VALUE(TEXT($E9,"DD"))),0),"")'

# ---------------------------------------------------------------------------
# 4) Resize the table / autofilter to cover the new rows.
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:K" + $endRow))

# ---------------------------------------------------------------------------
# 5) Extend the data validations to cover the new rows.
# ---------------------------------------------------------------------------
$ws.Range("F2:J31").Validation.Delete()
$ws.Range("B2:B31").Validation.Delete()

$vWhole = $ws.Range("F2:J" + $endRow).Validation
$vWhole.Add(1, 1, 1, "0", "1")
$vWhole.ErrorTitle = "Error de Tipeo"
$vWhole.ErrorMessage = "solo es 1 u 0 para definir v o f"

$vCustom = $ws.Range("B2:B" + $endRow).Validation
$vCustom.Add(7, 1, 1, "COUNTIF($B$2:$B$668,$B2)=1")
$vCustom.ErrorTitle = "Usuario Existente"
$vCustom.ErrorMessage = "corrija el usuario o verifique si son correctos los datos"

# ---------------------------------------------------------------------------
# 6) Update the active selection, matching what was captured when the file
#    was last saved.
# ---------------------------------------------------------------------------
$ws.Range("O10").Select()

Write-Host "done"
